$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.501409000000001
$ws.Range("M2").Value = 1.667434
$ws.Range("N2").Value = 5.002302
$ws.Range("O2").Value = 0.3223739883484499
$ws.Range("P2").Value = 0.32237398834845
$ws.Range("Q2").Value = 2.501934138168667
$ws.Range("R2").Value = 22.517407243518
$ws.Range("S2").Value = 0.3205636554758203
$ws.Range("T2").Value = 0.3205636554758204
$ws.Range("H3").Value = 4.501409000000001
$ws.Range("O3").Value = 0.2193354457157105
$ws.Range("P3").Value = 0.2193354457157106
$ws.Range("S3").Value = 0.2181037391206894
$ws.Range("T3").Value = 0.2181037391206896
$ws.Range("H4").Value = 4.501409000000001
$ws.Range("M4").Value = 0.2055123333333333
$ws.Range("N4").Value = 0.616537
$ws.Range("O4").Value = 0.03973280534729575
$ws.Range("P4").Value = 0.03973280534729576
$ws.Range("Q4").Value = 0.3083650222925556
$ws.Range("R4").Value = 2.775285200633
$ws.Range("S4").Value = 0.03950968063425515
$ws.Range("T4").Value = 0.03950968063425516
$ws.Range("H5").Value = 4.501409000000001
$ws.Range("M5").Value = 1.661741333333333
$ws.Range("N5").Value = 4.985224
$ws.Range("O5").Value = 0.3212733944672698
$ws.Range("P5").Value = 0.3212733944672699
$ws.Range("Q5").Value = 2.493392464512889
$ws.Range("R5").Value = 22.440532180616
$ws.Range("S5").Value = 0.3194692421220851
$ws.Range("T5").Value = 0.3194692421220852
$ws.Range("H6").Value = 4.501409000000001
$ws.Range("M6").Value = 0.1178836666666667
$ws.Range("N6").Value = 0.353651
$ws.Range("O6").Value = 0.02279108365576842
$ws.Range("P6").Value = 0.02279108365576842
$ws.Range("Q6").Value = 0.1768808660287778
$ws.Range("R6").Value = 1.591927794259
$ws.Range("S6").Value = 0.02266309737450464
$ws.Range("T6").Value = 0.02266309737450465
$ws.Range("H7").Value = 4.501409000000001
$ws.Range("M7").Value = 0.385306
$ws.Range("N7").Value = 1.155918
$ws.Range("O7").Value = 0.07449328246550557
$ws.Range("P7").Value = 0.0744932824655056
$ws.Range("Q7").Value = 0.5781399653846667
$ws.Range("R7").Value = 5.203259688462
$ws.Range("S7").Value = 0.07407495579241301
$ws.Range("T7").Value = 0.07407495579241304
$ws.Range("I8").Value = 0.005615629480232302
$ws.Range("J8").Value = 0.005615629480232303
$ws.Range("M8").Value = 1.667434
$ws.Range("N8").Value = 5.002302
$ws.Range("O8").Value = 0.3223739883484499
$ws.Range("P8").Value = 0.32237398834845
$ws.Range("Q8").Value = 0.01412927990466667
$ws.Range("R8").Value = 0.127163519142
$ws.Range("S8").Value = 0.00181033287262962
$ws.Range("T8").Value = 0.001810332872629621
$ws.Range("I9").Value = 0.005615629480232302
$ws.Range("J9").Value = 0.005615629480232303
$ws.Range("O9").Value = 0.2193354457157105
$ws.Range("P9").Value = 0.2193354457157106
$ws.Range("S9").Value = 0.001231706595021036
$ws.Range("T9").Value = 0.001231706595021036
$ws.Range("I10").Value = 0.005615629480232302
$ws.Range("J10").Value = 0.005615629480232303
$ws.Range("M10").Value = 0.2055123333333333
$ws.Range("N10").Value = 0.616537
$ws.Range("O10").Value = 0.03973280534729575
$ws.Range("P10").Value = 0.03973280534729576
$ws.Range("Q10").Value = 0.001741443008555555
$ws.Range("R10").Value = 0.015672987077
$ws.Range("S10").Value = 0.0002231247130406057
$ws.Range("T10").Value = 0.0002231247130406058
$ws.Range("I11").Value = 0.005615629480232302
$ws.Range("J11").Value = 0.005615629480232303
$ws.Range("M11").Value = 1.661741333333333
$ws.Range("N11").Value = 4.985224
$ws.Range("O11").Value = 0.3212733944672698
$ws.Range("P11").Value = 0.3212733944672699
$ws.Range("Q11").Value = 0.01408104214488889
$ws.Range("R11").Value = 0.126729379304
$ws.Range("S11").Value = 0.001804152345184702
$ws.Range("T11").Value = 0.001804152345184702
$ws.Range("I12").Value = 0.005615629480232302
$ws.Range("J12").Value = 0.005615629480232303
$ws.Range("M12").Value = 0.1178836666666667
$ws.Range("N12").Value = 0.353651
$ws.Range("O12").Value = 0.02279108365576842
$ws.Range("P12").Value = 0.02279108365576842
$ws.Range("Q12").Value = 0.0009989068967777776
$ws.Range("R12").Value = 0.008990162070999999
$ws.Range("S12").Value = 0.0001279862812637737
$ws.Range("T12").Value = 0.0001279862812637737
$ws.Range("I13").Value = 0.005615629480232302
$ws.Range("J13").Value = 0.005615629480232303
$ws.Range("M13").Value = 0.385306
$ws.Range("N13").Value = 1.155918
$ws.Range("O13").Value = 0.07449328246550557
$ws.Range("P13").Value = 0.0744932824655056
$ws.Range("Q13").Value = 0.003264954608666666
$ws.Range("R13").Value = 0.029384591478
$ws.Range("S13").Value = 0.0004183266730925651
$ws.Range("T13").Value = 0.0004183266730925654
